$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update text values (Description / Value-MPN columns) ---

# Row 2 (C1 - 0603 Capacitor): value "0.1µF" -> "0.1 µF"
$ws.Range("D2").Value = "0.1 µF"

# Row 4 (D1 - Indicator LED): value "LED" -> "Red"
$ws.Range("D4").Value = "Red"

# Row 5 (J1 - RJ45 Jack): value "RJ45_pyControl" -> "0855135013"
# Force text format so the leading zero of the part number is preserved.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0855135013"

# Row 6 (R1 - 0603 Resistor): value "10KΩ" -> "20KΩ"
$ws.Range("D6").Value = "20KΩ"

# Row 7 (U1): description "ATtiny 84 MCU" -> "ATtiny 24 MCU", value "ATtiny84-20SSU" -> "ATTINY24A-SSUR"
$ws.Range("C7").Value = "ATtiny 24 MCU"
$ws.Range("D7").Value = "ATTINY24A-SSUR"

# Row 8 (X1 - ISP Header): value "AVR_ISP" -> "75869-331LF"
$ws.Range("D8").Value = "75869-331LF"

# --- Rebuild the hyperlinks collection with the target order/targets ---
# (E2 and E4 are brand-new links, E7's target changed, E3/E5/E8 targets are
#  unchanged text but need to be regenerated so the whole set is in the
#  desired order.)
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("E2"), "https://www.digikey.com/products/en?keywords=1276-1258-1-ND")
$ws.Hyperlinks.Add($ws.Range("E3"), "https://www.digikey.com/products/en?keywords=P15094CT-ND")
$ws.Hyperlinks.Add($ws.Range("E4"), "https://www.digikey.com/products/en?keywords=160-1447-1-ND")
$ws.Hyperlinks.Add($ws.Range("E5"), "https://www.digikey.com/products/en?keywords=WM3553CT-ND")
$ws.Hyperlinks.Add($ws.Range("E7"), "https://www.digikey.com/products/en?keywords=ATTINY24A-SSURCT-ND")
$ws.Hyperlinks.Add($ws.Range("E8"), "https://www.digikey.com/products/en?keywords=609-5122-ND")

# Re-apply the built-in Hyperlink cell style to every linked cell so they all
# share the same style record (rather than the ad-hoc one Hyperlinks.Add
# creates on its own).
$ws.Range("E2").Style = "Hyperlink"
$ws.Range("E3").Style = "Hyperlink"
$ws.Range("E4").Style = "Hyperlink"
$ws.Range("E5").Style = "Hyperlink"
$ws.Range("E7").Style = "Hyperlink"
$ws.Range("E8").Style = "Hyperlink"

$wb.Save()
